$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original row 10 held the "Costo Totale Preventivo" total in G10/H10 (D10/E10 are
# just empty placeholder cells, untouched throughout). That total line moves down
# to row 11, a new "Budget" line takes over G10/H10, and a new "Budget Rimanente"
# line is appended as row 12.

# 1) Clone row 10's current formatting (style s="1": border + center horizontal)
#    onto row 11 BEFORE row 10's own look is changed.
$ws.Range("G10:H10").Copy()
$ws.Range("G11:H11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-create the "Costo Totale Preventivo" line on row 11.
$ws.Range("G11").Value = "Costo Totale Preventivo"
$ws.Range("H11").Formula = "=SUM(H3:H8)"

# 2) Re-purpose row 10 as the new "Budget" line: border + center-h + center-v
#    (same look already used elsewhere on the sheet, e.g. H8).
$ws.Range("H8").Copy()
$ws.Range("G10:H10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("G10").Value = "Budget "
$ws.Range("H10").Value = 300000

# 3) New row 12: "Budget Rimanente" = Budget (H10) - Costo Totale Preventivo (H11).
$ws.Range("H8").Copy()
$ws.Range("G12:H12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("G12").Value = "Budget Rimanente"
$ws.Range("H12").Formula = "=H10-H11"

$ws.Range("H13").Select()
